# Adapt the "Power - Links" sheet so it lists one Pmax row per node
# (Node_X -> Node_X.1) instead of the prior branch-to-branch link list,
# and trim the now-unused trailing rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New From/To node pairs for rows 8-16 (one row per node, self-referential
# "Node_X" -> "Node_X.1" pair used for the Pmax calculation).
$fromNodes = @("Node_1", "Node_2", "Node_3", "Node_4", "Node_5", "Node_6", "Node_7", "Node_8", "Node_9")
$toNodes   = @("Node_1.1", "Node_2.1", "Node_3.1", "Node_4.1", "Node_5.1", "Node_6.1", "Node_7.1", "Node_8.1", "Node_9.1")

for ($i = 0; $i -lt $fromNodes.Count; $i++) {
    $r = 8 + $i
    $ws.Cells.Item($r, 3).Value = $fromNodes[$i]
    $ws.Cells.Item($r, 4).Value = $toNodes[$i]
}

# The table used to have 13 data rows (8-20); it now only needs 9 (8-16),
# so remove the trailing four rows entirely (shrinks the used range/dimension).
$ws.Rows("17:20").Delete()

# Update the remembered selection/active cell on the sheet.
$ws.Range("I6").Select()
